# Generate Report for Handback
# - Updates the "Ready for handoff" status to "Handback transform failed"
#   for the 20ef5a27-... file (zh-cn and de-de locales), on the Overview
#   sheet as well as the locale-specific sheets.
# - Records the handback transform error detail message for both the
#   zh-cn and de-de locale sheets, and widens the "Error Detail" column
#   to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhError = "Handback file name: 0gwmdvxp.ao1 is different with handoff file name: 20ef5a27-bae2-4dd7-9058-e2d9faf21e32.63479792decaaa07924eaf5ab13f6127ce8a998d.zh-cn."
$deError  = "Handback file name: 0gwmdvxp.ao1 is different with handoff file name: 20ef5a27-bae2-4dd7-9058-e2d9faf21e32.63479792decaaa07924eaf5ab13f6127ce8a998d.de-de."

# --- Overview sheet: update status for the 20ef5a27-... row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: update status + error detail, widen column P ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = $zhError
$wsZhCn.Range("P1").ColumnWidth = 39.17

# --- de-de sheet: update status + error detail, widen column P ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = $deError
$wsDeDe.Range("P1").ColumnWidth = 39.17
